# May 28 2020 update
# Updates computed CASES/CHILDREN-related figures (columns E, F, G) for rows 3-12
# on the active worksheet of Table_1.xlsx, reflecting refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 915763.7357711862
$ws.Range("F3").Value = 163505.01289087508
$ws.Range("G3").Value = 17.9

$ws.Range("E4").Value = 171037.82777760012
$ws.Range("F4").Value = 16664.416194232548
$ws.Range("G4").Value = 9.7

$ws.Range("E5").Value = 86921.60302080003
$ws.Range("F5").Value = 11340.8056353941
$ws.Range("G5").Value = 13.0

$ws.Range("E6").Value = 657804.3049727976
$ws.Range("F6").Value = 135499.7910612484
$ws.Range("G6").Value = 20.6

$ws.Range("E7").Value = 158705.82012960006
$ws.Range("F7").Value = 29714.801234687355
$ws.Range("G7").Value = 18.7

$ws.Range("E8").Value = 230822.24114880018
$ws.Range("F8").Value = 39800.924230264754
$ws.Range("G8").Value = 17.2

$ws.Range("E9").Value = 272825.2078272002
$ws.Range("F9").Value = 46702.59660070997
$ws.Range("G9").Value = 17.1

$ws.Range("E10").Value = 32494.692268800005
$ws.Range("F10").Value = 6787.807567225531
$ws.Range("G10").Value = 20.9

$ws.Range("E11").Value = 220747.3237008001
$ws.Range("F11").Value = 40466.90267813935
$ws.Range("G11").Value = 18.3

$ws.Range("E12").Value = 168.450696
$ws.Range("F12").Value = 31.98057984809888
$ws.Range("G12").Value = 19.0
